# Briefing-Polaris.docx edit:
#   1. Remove the stray "_GoBack" bookmark sitting right after
#      "Internet (Wi-Fi);" (Word leaves one of these at the location of the
#      last edit made before the file was saved; it is not meaningful
#      content and disappears once a newer edit supersedes it).
#   2. Fix the typo'd deadline date "14/32/2024." -> "14/03/2024." and let
#      Word drop a fresh "_GoBack" bookmark at the spot of that edit (right
#      in the middle of the corrected day-of-month digits), which is what
#      naturally happens after Word's last save following an edit.

$d = $word.ActiveDocument

# --- 1. Drop the old "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Correct the deadline date -------------------------------------
# The paragraph currently reads "...inicial é 14/32/2024." (day digits
# swapped/garbled). Locate the unique "32/2024." tail so we can pin down
# exactly where the "3" that starts the day begins.
$find = $d.Content
$find.Find.ClearFormatting()
$found = $find.Find.Execute("32/2024.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the '32/2024.' deadline text to correct"
}
$dayStart = $find.Start

# Type the missing leading "0" of the day in front of the existing "3".
$d.Range($dayStart, $dayStart).InsertBefore("0")

# Word re-plants "_GoBack" at the caret position right after what was
# just typed, i.e. between the new "0" and the original "3".
$caret = $dayStart + 1
$d.Bookmarks.Add("_GoBack", $d.Range($caret, $caret))

# Remove the stray extra "2" that followed the original "3" (the day had
# been "32" instead of "03"), collapsing "3" + "2/2024." into "3/2024.".
$strayTwo = $dayStart + 2
$d.Range($strayTwo, $strayTwo + 1).Delete()
